$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Write the new date as literal text (matching the existing Date column,
# which stores dates as shared strings rather than real date values).
# A plain .Value assignment would be smart-parsed into a date serial, so
# stage it as a text formula result and paste-special the value in.
$ws.Cells.Item(1, 10).Formula = "=""2025-12-27"""
$ws.Cells.Item(1, 10).Copy()
$ws.Cells.Item(83, 1).PasteSpecial(-4163)
$ws.Cells.Item(1, 10).ClearContents()

$ws.Cells.Item(83, 2).Value = 0
$ws.Cells.Item(83, 3).Value = 30
